$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4126.857
$ws.Range("I51").Value = 1801
$ws.Range("J51").Value = 4305.769
$ws.Range("K51").Value = 1801
$ws.Range("L51").Value = 4305.769
$ws.Range("M51").Value = -1317
$ws.Range("N51").Value = -5273.769
$ws.Range("H98").Value = 2583.7576
$ws.Range("I98").Value = 1405.6552
$ws.Range("K98").Value = 1405.6552
$ws.Range("M98").Value = 92.34480000000008
$ws.Range("H122").Value = 2583.7576
$ws.Range("I122").Value = 1405.6552
$ws.Range("K122").Value = 4216.9656
$ws.Range("M122").Value = -1766.9656
$ws.Range("H125").Value = 1640.2858
$ws.Range("I125").Value = 1721.6666
$ws.Range("J125").Value = 1579.25
$ws.Range("K125").Value = 15494.9994
$ws.Range("L125").Value = 14213.25
$ws.Range("M125").Value = -13034.9994
$ws.Range("N125").Value = -19133.25
$ws.Range("H129").Value = 1036.5625
$ws.Range("J129").Value = 1071.1476
$ws.Range("L129").Value = 3213.4428
$ws.Range("N129").Value = -13213.4428
$ws.Range("H131").Value = 4715.5884
$ws.Range("I131").Value = 3253.889
$ws.Range("J131").Value = 6360
$ws.Range("K131").Value = 9761.667000000001
$ws.Range("L131").Value = 19080
$ws.Range("M131").Value = -4721.667000000001
$ws.Range("N131").Value = -29160
$ws.Range("H132").Value = 146774.72
$ws.Range("I132").Value = 184554.19
$ws.Range("J132").Value = 8250
$ws.Range("K132").Value = 553662.5700000001
$ws.Range("L132").Value = 24750
$ws.Range("M132").Value = -551132.5700000001
$ws.Range("N132").Value = -29810
$ws.Range("H137").Value = 4102.1577
$ws.Range("I137").Value = 3487.1614
$ws.Range("K137").Value = 10461.4842
$ws.Range("M137").Value = -7911.484199999999
$ws.Range("H141").Value = 2005.2188
$ws.Range("I141").Value = 1759.3572
$ws.Range("J141").Value = 3726.25
$ws.Range("K141").Value = 5278.071599999999
$ws.Range("L141").Value = 11178.75
$ws.Range("M141").Value = -98.07159999999931
$ws.Range("N141").Value = -21538.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 14940.167
$ws.Range("I41").Value = 2885.3333
$ws.Range("J41").Value = 26995
$ws.Range("K41").Value = 2885.3333
$ws.Range("L41").Value = 26995
$ws.Range("M41").Value = -2471.3333
$ws.Range("N41").Value = -27823
$ws.Range("H45").Value = 1349.2
$ws.Range("I45").Value = 1153
$ws.Range("K45").Value = 1153
$ws.Range("M45").Value = -776
$ws.Range("H95").Value = 49845
$ws.Range("J95").Value = 49845
$ws.Range("L95").Value = 49845
$ws.Range("N95").Value = -55337
$ws.Range("H110").Value = 688.65717
$ws.Range("I110").Value = 655.0741
$ws.Range("J110").Value = 802
$ws.Range("K110").Value = 655.0741
$ws.Range("L110").Value = 802
$ws.Range("M110").Value = 1389.9259
$ws.Range("N110").Value = -4892
$ws.Range("H122").Value = 2105.5557
$ws.Range("I122").Value = 1165.8572
$ws.Range("K122").Value = 3497.5716
$ws.Range("M122").Value = -1047.5716
$ws.Range("H132").Value = 2694.4187
$ws.Range("I132").Value = 1882.1
$ws.Range("J132").Value = 4569
$ws.Range("K132").Value = 5646.299999999999
$ws.Range("L132").Value = 13707
$ws.Range("M132").Value = -3116.299999999999
$ws.Range("N132").Value = -18767

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 901.13336
$ws.Range("I94").Value = 986.8570999999999
$ws.Range("K94").Value = 986.8570999999999
$ws.Range("M94").Value = -535.8570999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25004272
$ws.Range("I31").Value = 2210.182
$ws.Range("K31").Value = 2210.182
$ws.Range("M31").Value = -1915.182
$ws.Range("H34").Value = 25004272
$ws.Range("I34").Value = 2210.182
$ws.Range("K34").Value = 2210.182
$ws.Range("M34").Value = -2008.182
$ws.Range("H58").Value = 1920.6769
$ws.Range("I58").Value = 1658.1724
$ws.Range("J58").Value = 4095.7144
$ws.Range("K58").Value = 1658.1724
$ws.Range("L58").Value = 4095.7144
$ws.Range("M58").Value = -1455.1724
$ws.Range("N58").Value = -4501.7144
$ws.Range("H132").Value = 2659.275
$ws.Range("I132").Value = 1170.826
$ws.Range("J132").Value = 4673.0586
$ws.Range("K132").Value = 3512.478
$ws.Range("L132").Value = 14019.1758
$ws.Range("M132").Value = -982.4780000000001
$ws.Range("N132").Value = -19079.1758
$ws.Range("H134").Value = 11784.917
$ws.Range("I134").Value = 17903.166
$ws.Range("J134").Value = 5666.6665
$ws.Range("K134").Value = 53709.49800000001
$ws.Range("L134").Value = 16999.9995
$ws.Range("M134").Value = -51174.49800000001
$ws.Range("N134").Value = -22069.9995
$ws.Range("H136").Value = 1920.6769
$ws.Range("I136").Value = 1658.1724
$ws.Range("J136").Value = 4095.7144
$ws.Range("K136").Value = 4974.5172
$ws.Range("L136").Value = 12287.1432
$ws.Range("M136").Value = -2424.5172
$ws.Range("N136").Value = -17387.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 962.5946
$ws.Range("J131").Value = 856.8857400000001
$ws.Range("L131").Value = 2570.65722
$ws.Range("N131").Value = -12650.65722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3021.9
$ws.Range("I126").Value = 2970.9678
$ws.Range("J126").Value = 3698.5715
$ws.Range("K126").Value = 8912.903399999999
$ws.Range("L126").Value = 11095.7145
$ws.Range("M126").Value = -6442.903399999999
$ws.Range("N126").Value = -16035.7145
$ws.Range("H132").Value = 4222.375
$ws.Range("I132").Value = 2912
$ws.Range("J132").Value = 5008.6
$ws.Range("K132").Value = 8736
$ws.Range("L132").Value = 15025.8
$ws.Range("M132").Value = -6206
$ws.Range("N132").Value = -20085.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1918.3636
$ws.Range("J22").Value = 2516.75
$ws.Range("L22").Value = 2516.75
$ws.Range("N22").Value = -3106.75
$ws.Range("H27").Value = 1918.3636
$ws.Range("J27").Value = 2516.75
$ws.Range("L27").Value = 2516.75
$ws.Range("N27").Value = -2730.75
$ws.Range("H32").Value = 784.75
$ws.Range("I32").Value = 784.75
$ws.Range("K32").Value = 784.75
$ws.Range("M32").Value = -467.75
$ws.Range("H93").Value = 2024.3334
$ws.Range("I93").Value = 1223.5555
$ws.Range("J93").Value = 2825.111
$ws.Range("K93").Value = 1223.5555
$ws.Range("L93").Value = 2825.111
$ws.Range("M93").Value = 24.44450000000006
$ws.Range("N93").Value = -5321.111
$ws.Range("H132").Value = 4070.7302
$ws.Range("I132").Value = 1808.04
$ws.Range("K132").Value = 5424.12
$ws.Range("M132").Value = -2894.12

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1444.6296
$ws.Range("I136").Value = 661.5294
$ws.Range("J136").Value = 2775.9
$ws.Range("K136").Value = 1984.5882
$ws.Range("L136").Value = 8327.700000000001
$ws.Range("M136").Value = 565.4117999999999
$ws.Range("N136").Value = -13427.7
